$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "42.046.14"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.250.90"
$ws.Range("E3").Value = "  -0.94%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 (BNB)
$ws.Range("E5").Value = "  +0.23%  "

# Row 6 (Solana) - numeric-looking value, force text with leading apostrophe
$ws.Range("D6").Value = "'96.87"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  -1.30%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 (Cardano)
$ws.Range("E9").Value = "  -0.61%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'35.16"
$ws.Range("E10").Value = "  -1.73%  "

# Row 11 (Dogecoin)
$ws.Range("E11").Value = "  +2.42%  "

# Row 12 (TRON)
$ws.Range("E12").Value = "  +1.66%  "

# Row 13 (Polkadot)
$ws.Range("D13").Value = "'6.74"
$ws.Range("E13").Value = "  +1.12%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "2.600.48"
$ws.Range("E14").Value = "  -0.91%  "

# Row 15 (Chainlink)
$ws.Range("E15").Value = "  +0.46%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.247.17"
$ws.Range("E16").Value = "  -1.43%  "

# Row 17 (Polygon)
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  -2.06%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "42.054.61"

# Row 19 (InternetComputer(DFINITY))
$ws.Range("D19").Value = "'12.15"
$ws.Range("E19").Value = "  -3.01%  "

# Row 20 (ShibaInu)
$ws.Range("D20").Value = "0.0₃0904"

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  -0.49%  "

# Row 22 (Litecoin)
$ws.Range("D22").Value = "'67.07"
$ws.Range("E22").Value = "  -0.79%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'235.83"
$ws.Range("E23").Value = "  -1.97%  "

# Row 24 (PancakeSwap)
$ws.Range("D24").Value = "'2.57"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25 (ImmutableX)
$ws.Range("E25").Value = "  -1.11%  "

# Row 26 (Dai)
$ws.Range("D26").Value = "'0.990"
$ws.Range("E26").Value = "  -0.88%  "

# Row 27 (InjectiveProtocol)
$ws.Range("D27").Value = "'37.92"
$ws.Range("E27").Value = "  +0.97%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").Value = "'23.31"
$ws.Range("E28").Value = "  -2.14%  "

# Row 29 and Row 30 swap: Cosmos (was row29) <-> Toncoin (was row30), with updated data
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.50"
$ws.Range("E30").Value = "  -0.38%  "

# Row 31 (Monero)
$ws.Range("D31").Value = "'167.37"
$ws.Range("E31").Value = "  +4.65%  "

# Row 32 (FirstDigitalUSD)
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.08%  "

# Row 33 (Filecoin)
$ws.Range("D33").Value = "'5.16"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34 (Celestia)
$ws.Range("D34").Value = "'17.54"
$ws.Range("E34").Value = "  +2.52%  "

# Row 35 (LidoDAOToken)
$ws.Range("E35").Value = "  -3.84%  "

# Row 36 (Hedera)
$ws.Range("D36").Value = "'0.0719"
$ws.Range("E36").Value = "  -3.06%  "

# Row 37 (WEMIXToken)
$ws.Range("E37").Value = "  +1.65%  "

# Row 38 (Stellar)
$ws.Range("E38").Value = "  -0.01%  "

# Row 39 (Kaspa)
$ws.Range("E39").Value = "  -2.50%  "

# Row 40 (ARBITRUM)
$ws.Range("E40").Value = "  -2.40%  "

# Row 41 (RenderToken)
$ws.Range("D41").Value = "'4.07"
$ws.Range("E41").Value = "  -0.02%  "

# Row 42 (Maker)
$ws.Range("D42").Value = "1.939.50"
$ws.Range("E42").Value = "  -2.84%  "

# Row 43 (VeChain)
$ws.Range("E43").Value = "  -1.54%  "

# Row 44 (EnergySwap)
$ws.Range("D44").Value = "'18.58"
$ws.Range("E44").Value = "  -1.50%  "

# Row 45 (ApeXProtocol)
$ws.Range("E45").Value = "  -10.15%  "

# Row 46 (NEARProtocol)
$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = "  -1.86%  "

# Row 47 (FraxShare)
$ws.Range("D47").Value = "'9.67"
$ws.Range("E47").Value = "  -3.20%  "

# Row 48 (MultiversX)
$ws.Range("D48").Value = "'53.98"
$ws.Range("E48").Value = "  +1.31%  "

# Row 49 (RocketPoolETH)
$ws.Range("D49").Value = "2.469.93"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50 (BitcoinSV)
$ws.Range("D50").Value = "'71.09"
$ws.Range("E50").Value = "  -1.55%  "

# Row 51 (Aave)
$ws.Range("D51").Value = "'91.11"
$ws.Range("E51").Value = "  -0.51%  "
